$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update: the Odds (F) column no longer has fresh data for these
# fixtures, so clear the previously-populated values for rows 6-9
# (F6:F9), turning those numeric "1.84" cells into blank cells.
$ws.Range("F6:F9").ClearContents()
